# Auto-generated edit script applying the Sagittarius_Profits diff
# (workbook sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR correspond to the
#  concatenated per-sheet blocks in the source diff).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1235.1765
$ws.Range("I2").Value = 243.58333
$ws.Range("J2").Value = 3615
$ws.Range("K2").Value = 243.58333
$ws.Range("L2").Value = 3615
$ws.Range("M2").Value = -130.58333
$ws.Range("N2").Value = -3841
$ws.Range("H70").Value = 2211.111
$ws.Range("I70").Value = 425
$ws.Range("K70").Value = 1275
$ws.Range("M70").Value = -1005
$ws.Range("H73").Value = 2211.111
$ws.Range("I73").Value = 425
$ws.Range("K73").Value = 1275
$ws.Range("M73").Value = -339
$ws.Range("H106").Value = 31213.5
$ws.Range("I106").Value = 36756.2
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 36756.2
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -36125.2
$ws.Range("N106").Value = -4762
$ws.Range("H107").Value = 676.7
$ws.Range("J107").Value = 787.6
$ws.Range("L107").Value = 787.6
$ws.Range("N107").Value = -4627.6
$ws.Range("H135").Value = 188.28572
$ws.Range("I135").Value = 216.27272
$ws.Range("J135").Value = 85.666664
$ws.Range("K135").Value = 1946.45448
$ws.Range("L135").Value = 770.9999759999999
$ws.Range("M135").Value = 588.5455200000001
$ws.Range("N135").Value = -5840.999976
$ws.Range("H137").Value = 1303
$ws.Range("I137").Value = 1303
$ws.Range("K137").Value = 3909
$ws.Range("M137").Value = -1359
$ws.Range("H138").Value = 4394.484
$ws.Range("J138").Value = 5876.8887
$ws.Range("L138").Value = 17630.6661
$ws.Range("N138").Value = -27910.6661

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2911.5557
$ws.Range("I32").Value = 2994.0667
$ws.Range("K32").Value = 2994.0667
$ws.Range("M32").Value = -2707.0667
$ws.Range("H61").Value = 1934.8889
$ws.Range("I61").Value = 1996.9412
$ws.Range("K61").Value = 1996.9412
$ws.Range("M61").Value = -1784.9412
$ws.Range("H136").Value = 1934.8889
$ws.Range("I136").Value = 1996.9412
$ws.Range("K136").Value = 5990.8236
$ws.Range("M136").Value = -3440.8236

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 892.46155
$ws.Range("I94").Value = 719.875
$ws.Range("K94").Value = 719.875
$ws.Range("M94").Value = -268.875
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 93362.46000000001
$ws.Range("I94").Value = 163815.28
$ws.Range("J94").Value = 11167.5
$ws.Range("K94").Value = 163815.28
$ws.Range("L94").Value = 11167.5
$ws.Range("M94").Value = -163364.28
$ws.Range("N94").Value = -12069.5
$ws.Range("H99").Value = 1791.2307
$ws.Range("I99").Value = 1411.375
$ws.Range("J99").Value = 2399
$ws.Range("K99").Value = 1411.375
$ws.Range("L99").Value = 2399
$ws.Range("M99").Value = 86.625
$ws.Range("N99").Value = -5395
$ws.Range("H105").Value = 2645.7058
$ws.Range("I105").Value = 1297.4445
$ws.Range("K105").Value = 1297.4445
$ws.Range("M105").Value = 449.5554999999999
$ws.Range("H126").Value = 1791.2307
$ws.Range("I126").Value = 1411.375
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 4234.125
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -1764.125
$ws.Range("N126").Value = -12137
$ws.Range("H130").Value = 84798.8
$ws.Range("J130").Value = 84798.8
$ws.Range("L130").Value = 84798.8
$ws.Range("N130").Value = -94838.8
$ws.Range("H134").Value = 828.5
$ws.Range("I134").Value = 828.5
$ws.Range("K134").Value = 2485.5
$ws.Range("M134").Value = 49.5
$ws.Range("H138").Value = 65756
$ws.Range("J138").Value = 65756
$ws.Range("L138").Value = 65756
$ws.Range("N138").Value = -76036

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 835.9524
$ws.Range("I12").Value = 875.75
$ws.Range("J12").Value = 782.8889
$ws.Range("K12").Value = 2627.25
$ws.Range("L12").Value = 2348.6667
$ws.Range("M12").Value = -2454.25
$ws.Range("N12").Value = -2694.6667
$ws.Range("H138").Value = 7144455.5
$ws.Range("I138").Value = 9092535
$ws.Range("K138").Value = 27277605
$ws.Range("M138").Value = -27272465

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7000
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 4000
$ws.Range("N55").Value = -4654
$ws.Range("H70").Value = 7999.5
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730
$ws.Range("H73").Value = 7999.5
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064
$ws.Range("H80").Value = 3606.4285
$ws.Range("J80").Value = 499
$ws.Range("L80").Value = 499
$ws.Range("N80").Value = -2495
$ws.Range("H83").Value = 3606.4285
$ws.Range("J83").Value = 499
$ws.Range("L83").Value = 2495
$ws.Range("N83").Value = -12479
$ws.Range("H97").Value = 852.4286
$ws.Range("I97").Value = 1244.75
$ws.Range("J97").Value = 329.33334
$ws.Range("K97").Value = 1244.75
$ws.Range("L97").Value = 329.33334
$ws.Range("M97").Value = -748.75
$ws.Range("N97").Value = -1321.33334
$ws.Range("H102").Value = 3571.0908
$ws.Range("I102").Value = 3507.5
$ws.Range("J102").Value = 4207
$ws.Range("K102").Value = 3507.5
$ws.Range("L102").Value = 4207
$ws.Range("M102").Value = -1885.5
$ws.Range("N102").Value = -7451
$ws.Range("H107").Value = 1652.6786
$ws.Range("I107").Value = 1046.3125
$ws.Range("K107").Value = 1046.3125
$ws.Range("M107").Value = 873.6875
$ws.Range("H126").Value = 4374.8335
$ws.Range("I126").Value = 4369.8
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 13109.4
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -10639.4
$ws.Range("N126").Value = -18140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8664.866
$ws.Range("I7").Value = 5992
$ws.Range("J7").Value = 9076.076999999999
$ws.Range("K7").Value = 5992
$ws.Range("L7").Value = 9076.076999999999
$ws.Range("M7").Value = -5880
$ws.Range("N7").Value = -9300.076999999999
$ws.Range("H46").Value = 4020
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 5766.6665
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 5766.6665
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -6142.6665
$ws.Range("H61").Value = 4335
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3798
$ws.Range("H93").Value = 1739.3334
$ws.Range("I93").Value = 1329.8
$ws.Range("J93").Value = 2251.25
$ws.Range("K93").Value = 1329.8
$ws.Range("L93").Value = 2251.25
$ws.Range("M93").Value = -81.79999999999995
$ws.Range("N93").Value = -4747.25
$ws.Range("H113").Value = 4335
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -1830
$ws.Range("H126").Value = 8664.866
$ws.Range("I126").Value = 5992
$ws.Range("J126").Value = 9076.076999999999
$ws.Range("K126").Value = 17976
$ws.Range("L126").Value = 27228.231
$ws.Range("M126").Value = -15506
$ws.Range("N126").Value = -32168.231
$ws.Range("H132").Value = 5322.9287
$ws.Range("I132").Value = 5946.364
$ws.Range("K132").Value = 17839.092
$ws.Range("M132").Value = -15309.092

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H107").Value = 659.2105
$ws.Range("I107").Value = 614
$ws.Range("K107").Value = 1842
$ws.Range("M107").Value = 78
$ws.Range("H113").Value = 650.5
$ws.Range("I113").Value = 585.375
$ws.Range("J113").Value = 780.75
$ws.Range("K113").Value = 1756.125
$ws.Range("L113").Value = 2342.25
$ws.Range("M113").Value = 413.875
$ws.Range("N113").Value = -6682.25
$ws.Range("H136").Value = 3950.5557
$ws.Range("I136").Value = 4109.609
$ws.Range("K136").Value = 12328.827
$ws.Range("M136").Value = -9778.827000000001
